# Apply cryptos list update (Sun Oct  1 16:00:50 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new value is a plain number string (e.g. "215.84") ---
# These must stay TEXT (matching original inlineStr cells), so we briefly
# force a Text number format while assigning, then restore the default
# "Normal" style so no stray formatting is left behind.
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "215.84"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "23.14"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.261"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "4.20"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "0.556"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "66.97"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "236.28"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "8.02"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "9.62"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "147.39"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "7.34"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "16.47"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "0.0506"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "0.948"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "69.19"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "5.75"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "90.19"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "8.26"
$c.NumberFormat = "General"
$c.Style = "Normal"

# --- Remaining cells: plain text values, safe to assign directly ---
$ws.Cells.Item(2, 4).Value = "27.192.99"
$ws.Cells.Item(2, 5).Value = "  +0.64%  "
$ws.Cells.Item(3, 4).Value = "1.685.99"
$ws.Cells.Item(3, 5).Value = "  +0.14%  "
$ws.Cells.Item(4, 5).Value = "  +0.05%  "
$ws.Cells.Item(5, 5).Value = "  -0.13%  "
$ws.Cells.Item(6, 5).Value = "  +0.48%  "
$ws.Cells.Item(7, 5).Value = "  +0.02%  "
$ws.Cells.Item(8, 5).Value = "  +7.86%  "
$ws.Cells.Item(9, 5).Value = "  +2.99%  "
$ws.Cells.Item(10, 5).Value = "  +0.82%  "
$ws.Cells.Item(11, 5).Value = "  +0.28%  "
$ws.Cells.Item(12, 4).Value = "1.923.01"
$ws.Cells.Item(13, 4).Value = "1.685.59"
$ws.Cells.Item(13, 5).Value = "  -0.25%  "
$ws.Cells.Item(14, 5).Value = "  +2.24%  "
$ws.Cells.Item(15, 5).Value = "  +3.92%  "
$ws.Cells.Item(16, 5).Value = "  +1.04%  "
$ws.Cells.Item(17, 4).Value = "27.188.57"
$ws.Cells.Item(17, 5).Value = "  +0.43%  "
$ws.Cells.Item(18, 5).Value = "  -0.46%  "
$ws.Cells.Item(19, 5).Value = "  -2.71%  "
$ws.Cells.Item(20, 5).Value = "  +1.20%  "
$ws.Cells.Item(21, 5).Value = "  +0.02%  "
$ws.Cells.Item(22, 5).Value = "  +2.25%  "
$ws.Cells.Item(23, 5).Value = "  +3.93%  "
$ws.Cells.Item(24, 5).Value = "  -2.42%  "
$ws.Cells.Item(25, 5).Value = "  +0.37%  "
$ws.Cells.Item(26, 5).Value = "  +1.32%  "
$ws.Cells.Item(27, 5).Value = "  +2.38%  "
$ws.Cells.Item(28, 5).Value = "  +0.57%  "
$ws.Cells.Item(29, 5).Value = "  +0.11%  "
$ws.Cells.Item(30, 5).Value = "  +1.26%  "
$ws.Cells.Item(31, 5).Value = "  +0.10%  "
$ws.Cells.Item(32, 5).Value = "  +1.29%  "
$ws.Cells.Item(33, 4).Value = "1.539.80"
$ws.Cells.Item(34, 5).Value = "  +1.96%  "
$ws.Cells.Item(35, 5).Value = "  -1.29%  "
$ws.Cells.Item(36, 5).Value = "  +2.62%  "
$ws.Cells.Item(37, 5).Value = "  +3.23%  "
$ws.Cells.Item(38, 5).Value = "  -0.46%  "
$ws.Cells.Item(39, 5).Value = "  -0.27%  "
$ws.Cells.Item(40, 5).Value = "  +1.57%  "
$ws.Cells.Item(41, 2).Value = "Aave"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(41, 5).Value = "  +1.04%  "
$ws.Cells.Item(42, 2).Value = "FraxShare"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(42, 5).Value = "  -0.57%  "
$ws.Cells.Item(43, 5).Value = "  +0.03%  "
$ws.Cells.Item(44, 5).Value = "  -1.35%  "
$ws.Cells.Item(45, 4).Value = "1.830.52"
$ws.Cells.Item(45, 5).Value = "  +0.17%  "
$ws.Cells.Item(46, 5).Value = "  +0.85%  "
$ws.Cells.Item(47, 5).Value = "  -0.19%  "
$ws.Cells.Item(48, 5).Value = "  +17.16%  "
$ws.Cells.Item(49, 5).Value = "  +5.41%  "
$ws.Cells.Item(50, 5).Value = "  +4.11%  "
$ws.Cells.Item(51, 5).Value = "  -0.75%  "
